# Apply the "Clusters" sheet rework: the old lookup table (department-id -> name)
# that lived at AA2:AB7 is relocated to AC2:AD7, each row grows three new helper
# columns (Y = min distance, Z = matching cluster index, AA = cluster name from the
# min distance), and the existing W-column VLOOKUP is repointed at the relocated
# table. Finally the viewport/selection/column-width cosmetics are nudged to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clusters")

# Data rows present on the sheet (rows 9, 17 and 34 are blank separator rows).
$dataRows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,35,36,37,38,39,40,41,42,43,44,45,46,47)

# 1) Relocate the lookup table from AA2:AB7 to AC2:AD7.
$lookupVals = $ws.Range("AA2:AB7").Value2
$ws.Range("AC2:AD7").Value2 = $lookupVals
$ws.Range("AA2:AB7").ClearContents()

# 2) Re-point the existing W-column VLOOKUP at the relocated table, and add the
#    new Y (min distance), Z (matching index) and AA (cluster name) helper columns.
foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 23).Formula = "=VLOOKUP(V$r,AC`$2:AD`$7,2,FALSE)"        # W
    $ws.Cells.Item($r, 25).Formula = "=MIN(M" + $r + ":R" + $r + ")"            # Y
    $ws.Cells.Item($r, 26).Formula = "=MATCH(Y" + $r + ",M" + $r + ":R" + $r + ",FALSE)"  # Z
    $ws.Cells.Item($r, 27).Formula = "=VLOOKUP(Z$r,AC`$2:AD`$7,2,FALSE)"        # AA
}

# 3) New narrow columns O:S (the DISTSQ/MINDIST helper block).
$ws.Range("O1:S1").ColumnWidth = 2

# 4) Viewport + selection cosmetics.
$ws.Activate()
$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L1:L1048576").Select()

# 5) Workbook window position cosmetics.
$excel.Left = 16560
$excel.Top = 9150
$excel.Width = 20910
$excel.Height = 11835
